$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("第一学年")

# Fix category label typo: "其他" -> "其它" (matches existing "其它" entries)
$ws.Range("F23").Value = "其它"
$ws.Range("F32").Value = "其它"
$ws.Range("F39").Value = "其它"

# Add new expense entry (row 45): 2018-07-10, 支出生活费100
$ws.Range("B44:G44").Copy()
$ws.Range("B45").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B45").Value = 43
$ws.Range("C45").Value = "支出"
$ws.Range("D45").Value = 100
$ws.Range("E45").Value = (Get-Date -Year 2018 -Month 7 -Day 10).Date
$ws.Range("F45").Value = "生活费"
$ws.Range("G45").Value = "生活费(7/11- )"

# Extend the autofilter range to include the new row
$ws.AutoFilterMode = $false
$ws.Range("B2:G45").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=第一学年!`$B`$2:`$G`$45"
    }
}

# Recalculate formulas so cached values reflect the new data
$excel.CalculateFull()

# Reset the view: remove the scrolled-to-row-37 position, select K18
$excel.Goto($ws.Range("A1"))
$ws.Range("K18").Select()
